{"js": "const body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1) Move the \"_GoBack\" bookmark: it used to sit right before the final\n//    \".\" after \"(keyboard shortcuts vary depending on the OS)\"; it now\n//    sits right after the report date's day-of-month number.\n//    Remove the old one first (bookmark names must stay unique).\n// ---------------------------------------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Change the report date from \"May 11, 2017\" to \"May 14, 2017\".\n//    Scope the search for \"11\" to the date phrase so we don't touch the\n//    \"11\" that appears elsewhere in the document (e.g. in an IP address).\n// ---------------------------------------------------------------------\nconst dateMatches = body.search(\"May 11, 2017\", { matchCase: true });\ndateMatches.load(\"items\");\nawait context.sync();\n\nconst dateRange = dateMatches.items[0];\nconst dayMatches = dateRange.search(\"11\", { matchCase: true });\ndayMatches.load(\"items\");\nawait context.sync();\n\nconst dayRange = dayMatches.items[0];\ndayRange.insertText(\"14\", \"Replace\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 3) Insert the \"_GoBack\" bookmark right after the new \"14\", before the\n//    \", 2017\" that follows it.\n// ---------------------------------------------------------------------\nconst newDateMatches = body.search(\"May 14, 2017\", { matchCase: true });\nnewDateMatches.load(\"items\");\nawait context.sync();\n\nconst newDateRange = newDateMatches.items[0];\nconst newDayMatches = newDateRange.search(\"14\", { matchCase: true });\nnewDayMatches.load(\"items\");\nawait context.sync();\n\nconst newDayRange = newDayMatches.items[0];\nconst afterDayRange = newDayRange.getRange(\"End\");\nafterDayRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 4) The two placeholder \"Test\" bullet items (one under \"SQL Database\",\n//    one under \"Microsoft Azure\") become links to the Azure homepage.\n// ---------------------------------------------------------------------\nconst testMatches = body.search(\"Test\", { matchCase: true, matchWholeWord: true });\ntestMatches.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < testMatches.items.length; i++) {\n  const testRange = testMatches.items[i];\n  testRange.insertText(\"https://azure.microsoft.com/en-us/\", \"Replace\");\n  await context.sync();\n  testRange.hyperlink = \"https://azure.microsoft.com/en-us/\";\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# -----------------------------------------------------------------------\n# 1) Move the \"_GoBack\" bookmark: it used to sit right before the final\n#    \".\" after \"(keyboard shortcuts vary depending on the OS)\"; it now\n#    sits right after the report date's day-of-month number. Remove the\n#    old one first (bookmark names must stay unique).\n# -----------------------------------------------------------------------\n$goBack = $d.Bookmarks.Item(\"_GoBack\")\n$goBack.Delete()\n\n# -----------------------------------------------------------------------\n# 2) Change the report date from \"May 11, 2017\" to \"May 14, 2017\".\n#    Scope the Find to just the \"11\" so the \"May \" and \", 2017\" runs\n#    around it are left alone.\n# -----------------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"11\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\n$dayRange = $d.Range($find.Parent.Start, $find.Parent.End)\n$dayRange.Text = \"14\"\n\n# Replacing the text re-merges the three date runs (\"May \", the day\n# number, and \", 2017\") into a single run because they all share the\n# same visible formatting. Re-split \"May \" away from \"14\" by dropping a\n# bookmark right after \"May \" and immediately deleting it again -- the\n# bookmark insert/delete forces Word to split the run at that position\n# without leaving anything behind.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"May 14\"\n$find2.MatchCase = $true\n$found2 = $find2.Execute()\n$mayDayRange = $d.Range($find2.Parent.Start, $find2.Parent.End)\n$splitPoint = $mayDayRange.Start + 4\n$splitRange = $d.Range($splitPoint, $splitPoint)\n$d.Bookmarks.Add(\"TempSplit\", $splitRange)\n$d.Bookmarks.Item(\"TempSplit\").Delete()\n\n# -----------------------------------------------------------------------\n# 3) Insert the \"_GoBack\" bookmark right after the new \"14\", before the\n#    \", 2017\" that follows it.\n# -----------------------------------------------------------------------\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Text = \"May 14\"\n$find3.MatchCase = $true\n$found3 = $find3.Execute()\n$mayDayRange2 = $d.Range($find3.Parent.Start, $find3.Parent.End)\n$afterDay = $d.Range($mayDayRange2.End, $mayDayRange2.End)\n$d.Bookmarks.Add(\"_GoBack\", $afterDay)\n\n# -----------------------------------------------------------------------\n# 4) The two placeholder \"Test\" bullet items (one under \"SQL Database\",\n#    one under \"Microsoft Azure\") become links to the Azure homepage.\n# -----------------------------------------------------------------------\n$testFind1 = $d.Content.Find\n$testFind1.ClearFormatting()\n$testFind1.Text = \"Test\"\n$testFind1.MatchWholeWord = $true\n$testFind1.MatchCase = $true\n$testFound1 = $testFind1.Execute()\n$testRange1 = $d.Range($testFind1.Parent.Start, $testFind1.Parent.End)\n$d.Hyperlinks.Add($testRange1, \"https://azure.microsoft.com/en-us/\", [Type]::Missing, [Type]::Missing, \"https://azure.microsoft.com/en-us/\")\n\n$testFind2 = $d.Content.Find\n$testFind2.ClearFormatting()\n$testFind2.Text = \"Test\"\n$testFind2.MatchWholeWord = $true\n$testFind2.MatchCase = $true\n$testFound2 = $testFind2.Execute()\n$testRange2 = $d.Range($testFind2.Parent.Start, $testFind2.Parent.End)\n$d.Hyperlinks.Add($testRange2, \"https://azure.microsoft.com/en-us/\", [Type]::Missing, [Type]::Missing, \"https://azure.microsoft.com/en-us/\")\n"}
